# Add a new "GlobalsqaDemoSite" worksheet, positioned right after
# "DashboardPage" (i.e. before "PIMPage"), and populate it with the
# locator table for the Globalsqa demo site. This mirrors the commit
# "Added Locators for GlobalsqaDemoSite".

$wb = $excel.ActiveWorkbook

$after  = $wb.Worksheets.Item("DashboardPage")
$newWs  = $wb.Worksheets.Add($null, $after)
$newWs.Name = "GlobalsqaDemoSite"

# Fill column A (labels) and C (xpath values) in the same interleaved
# order the strings were originally authored in, then column B (the
# constant "XPATH" marker, already present in the shared-string table).
$newWs.Range("A1").Value = "dragAndDropButton"
$newWs.Range("A2").Value = "sourceImageLocator"
$newWs.Range("A3").Value = "destiantionTrashLocator"
$newWs.Range("C1").Value = "//a[text()='DragAndDrop']"
$newWs.Range("A4").Value = "closeAdPopUpLocator"
$newWs.Range("C4").Value = "//span[text()='Close']"
$newWs.Range("C3").Value = "//div[@id='trash']"
$newWs.Range("A5").Value = "frameLocator"
$newWs.Range("C5").Value = "//iframe[@class='demo-frame lazyloaded']"
$newWs.Range("C2").Value = "//h5[text()='High Tatras']/parent::li"

$newWs.Range("B1").Value = "XPATH"
$newWs.Range("B2").Value = "XPATH"
$newWs.Range("B3").Value = "XPATH"
$newWs.Range("B4").Value = "XPATH"
$newWs.Range("B5").Value = "XPATH"

# Column widths matching the source workbook (closest value the
# character-width/pixel grid allows).
$newWs.Columns.Item(1).ColumnWidth = 25.67
$newWs.Columns.Item(3).ColumnWidth = 43.67

# The new sheet becomes the active / selected tab, with D13 selected.
$newWs.Range("D13").Select() | Out-Null
